$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = "6_211115_1297_rf_with_3in1_no_valid"
$ws.Range("B7").Value = 0.56465
$ws.Range("C7").Value = "random forest, in 3in1 data set, with no valid set"

$ws.Range("B8").Select()
